$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.070637
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.22405966666667
$ws.Range("N2").Value = 84.672179
$ws.Range("O2").Value = 0.3816548478108986
$ws.Range("P2").Value = 0.3816548478108986
$ws.Range("Q2").Value = 0.6645543008914445
$ws.Range("R2").Value = 5.980988708023
$ws.Range("S2").Value = 0.001074570364822368
$ws.Range("T2").Value = 0.001074570364822368

$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.070637
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("N3").Value = 59.306181
$ws.Range("O3").Value = 0.2673191094302723
$ws.Range("P3").Value = 0.2673191094302723
$ws.Range("Q3").Value = 0.4654678563663334
$ws.Range("R3").Value = 4.189210707297
$ws.Range("S3").Value = 0.0007526517600709367
$ws.Range("T3").Value = 0.0007526517600709365

$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.070637
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("M4").Value = 25.95900466666667
$ws.Range("N4").Value = 77.877014
$ws.Range("O4").Value = 0.351026042758829
$ws.Range("P4").Value = 0.351026042758829
$ws.Range("Q4").Value = 0.6112220708797779
$ws.Range("R4").Value = 5.500998637918
$ws.Range("S4").Value = 0.0009883332675926137
$ws.Range("T4").Value = 0.0009883332675926135

$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 28.22405966666667
$ws.Range("N5").Value = 84.672179
$ws.Range("O5").Value = 0.3816548478108986
$ws.Range("P5").Value = 0.3816548478108986
$ws.Range("Q5").Value = 232.9247442511994
$ws.Range("R5").Value = 2096.322698260795
$ws.Range("S5").Value = 0.3766344256149114
$ws.Range("T5").Value = 0.3766344256149114

$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("N6").Value = 59.306181
$ws.Range("O6").Value = 0.2673191094302723
$ws.Range("P6").Value = 0.2673191094302723
$ws.Range("S6").Value = 0.2638026997787428
$ws.Range("T6").Value = 0.2638026997787428

$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("M7").Value = 25.95900466666667
$ws.Range("N7").Value = 77.877014
$ws.Range("O7").Value = 0.351026042758829
$ws.Range("P7").Value = 0.351026042758829
$ws.Range("S7").Value = 0.3464085226446625
$ws.Range("T7").Value = 0.3464085226446624

$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 28.22405966666667
$ws.Range("N8").Value = 84.672179
$ws.Range("O8").Value = 0.3816548478108986
$ws.Range("P8").Value = 0.3816548478108986
$ws.Range("Q8").Value = 2.440261606799889
$ws.Range("R8").Value = 21.962354461199
$ws.Range("S8").Value = 0.003945851831164838
$ws.Range("T8").Value = 0.003945851831164838

$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("N9").Value = 59.306181
$ws.Range("O9").Value = 0.2673191094302723
$ws.Range("P9").Value = 0.2673191094302723
$ws.Range("Q9").Value = 1.709210725995666
$ws.Range("S9").Value = 0.002763757891458578
$ws.Range("T9").Value = 0.002763757891458578

$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("M10").Value = 25.95900466666667
$ws.Range("N10").Value = 77.877014
$ws.Range("O10").Value = 0.351026042758829
$ws.Range("P10").Value = 0.351026042758829
$ws.Range("R10").Value = 20.199817768334
$ws.Range("S10").Value = 0.003629186846573887
$ws.Range("T10").Value = 0.003629186846573887
